# Splits the paragraph "Merkið vélmennið þannig að vel sjáist hver á" into
# two paragraphs: the original text stays in its own paragraph, and a new
# list paragraph is added containing the new instructions, with the
# pre-existing "_GoBack" bookmark relocated between the two new runs.

$d = $word.ActiveDocument

# Locate the paragraph whose text ends with "hver á" and collapse the
# range to right after that text (i.e. right before the existing
# bookmark / paragraph mark).
$rng = $d.Content
$found = $rng.Find.Execute("Merkið vélmennið þannig að vel sjáist hver á", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph text"
}

$rng.Collapse(0)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="is-IS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="is-IS"/></w:rPr><w:t xml:space="preserve">Takið mynd af vélmenni frá hlið framan og að ofan skilið síðan slóðinni í Innu </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="is-IS"/></w:rPr><w:t>verkefni 1</w:t></w:r></w:p>'

# InsertXML on a collapsed range replaces the paragraph the range is
# located in, so re-emit the (unchanged) first paragraph followed by the
# new second paragraph.
$closingParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007F4EE0" w:rsidRDefault="007F4EE0" w:rsidP="005A7E46"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="is-IS"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="is-IS"/></w:rPr><w:t>Merkið vélmennið þannig að vel sjáist hver á</w:t></w:r></w:p>'

$null = $rng.InsertXML($closingParagraphXml + $newParagraphXml)
